# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list to use concise,
# impact-focused accomplishment statements (4 bullets instead of 6).
#
# Original 6 bullets under the "Impact" Heading3:
#   1. Discovered systematic race coding errors ... 23% to 64%
#   2. Achieved 87% prediction accuracy for voter turnout ... ±2.1%
#   3. Built cloud-based data warehouse solutions on AWS ... 99.94% accuracy
#   4. Built redistricting platform used by thousands ... 89 organizations
#   5. Developed longitudinal data analysis methods ... response quality
#   6. Designed ETL pipelines using PySpark, dbt ... geospatial datasets
#
# New 4 bullets:
#   1. Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard
#   2. Reduced polling margins from ±4.2% to ±2.1%
#   3. Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%
#   4. Reduced polling costs while increasing quality
#
# Bullets 1-3 get their text replaced in place; bullets 4 and 5 are deleted
# outright (paragraph and all), and bullet 6's text becomes the new 4th bullet.

$d = $word.ActiveDocument

# Locate the "Impact" Heading3 paragraph inside "KEY ACHIEVEMENTS AND IMPACT"
# unambiguously, then work relative to it so we don't depend on hard-coded
# absolute paragraph indices (robust against any earlier/later section size
# differences).
$sectionHeadingText = "KEY ACHIEVEMENTS AND IMPACT"
$impactHeadingText = "Impact"

$total = $d.Paragraphs.Count
$sectionIdx = 0
for ($i = 1; $i -le $total; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    if ($txt -eq $sectionHeadingText) {
        $sectionIdx = $i
        break
    }
}
if ($sectionIdx -eq 0) {
    throw "Could not find section heading '$sectionHeadingText'"
}

$impactIdx = 0
for ($i = $sectionIdx + 1; $i -le $total; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    $txt = $txt.TrimEnd([char]13, [char]7)
    if ($txt -eq $impactHeadingText) {
        $impactIdx = $i
        break
    }
}
if ($impactIdx -eq 0) {
    throw "Could not find 'Impact' sub-heading after section heading"
}

# The six bullet paragraphs immediately follow the "Impact" heading.
$b1 = $impactIdx + 1
$b2 = $impactIdx + 2
$b3 = $impactIdx + 3
$b4 = $impactIdx + 4
$b5 = $impactIdx + 5
$b6 = $impactIdx + 6

# Sanity-check the bullets are what we expect before mutating anything.
$chk1 = $d.Paragraphs.Item($b1).Range.Text
$chk6 = $d.Paragraphs.Item($b6).Range.Text
if ($chk1 -notmatch "Discovered systematic race coding errors") {
    throw "Unexpected bullet 1 text: $chk1"
}
if ($chk6 -notmatch "Designed ETL pipelines") {
    throw "Unexpected bullet 6 text: $chk6"
}

# Replace text of bullets 1-3 in place (paragraph count unchanged).
$d.Paragraphs.Item($b1).Range.Text = "• Predictive excellence: Achieved 87% voter turnout accuracy vs. 71% industry standard"
$d.Paragraphs.Item($b2).Range.Text = "• Reduced polling margins from ±4.2% to ±2.1%"
$d.Paragraphs.Item($b3).Range.Text = "• Methodological advancement: Improved segmentation accuracy 34% and survey incidence 28%"

# Bullet 6 becomes the new bullet 4.
$d.Paragraphs.Item($b6).Range.Text = "• Reduced polling costs while increasing quality"

# Delete old bullets 4 and 5 entirely (paragraph mark included). Delete from
# the higher index down so earlier indices stay valid.
$d.Paragraphs.Item($b5).Range.Delete()
$d.Paragraphs.Item($b4).Range.Delete()

Write-Output "Done rewriting Key Achievements bullets."
